$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) cells per the scheduled data refresh.
# Values are written with a leading apostrophe so Excel stores them as literal text
# (matching the workbooks existing inlineStr/text cells), then the style is reset to
# Normal so the transient quote-prefix number format left behind by COM does not stick
# to the cell.
$ws.Range("D2").Value = "29.389.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "1.848.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'240.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6275"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2905"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'5.025"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.6784"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D16").Value = "'6.154"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "29.402.73"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'226.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D20").Value = "'0.9997"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'7.511"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'158.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.1380"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'8.399"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'17.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'1.384"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.463"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.05608"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.081"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.837"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.6942"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.580"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "1.232.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.01803"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.719"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'6.396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.9068"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'101.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'66.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'7.175"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00000000118"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.4012"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'8.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.678"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.1141"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.05703"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = "'0.4625"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
$ws.Range("E51").Style = "Normal"
